# Regenerate merged AHB file:
#  - header row relabelled from the generic "_old"/"_new" suffixes to the
#    concrete format versions "_FV2210"/"_FV2304"
#  - data range turned into a proper Excel Table ("Table1") with AutoFilter
#  - header row frozen so it stays visible while scrolling

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the 21 header cells (A1:U1).
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# 2. Convert A1:U75 into an Excel Table ("Table1") with a header row / AutoFilter.
$dataRange = $ws.Range("A1:U75")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# 3. Freeze the header row (split below row 1, so row 1 stays visible).
[void]$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Regenerated merged AHB file: headers relabelled, Table1 created, header row frozen."
